$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 76, shifting existing rows 76..148 down to 77..149
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new weekly price entry
$ws.Range("A76").Value = 1
$ws.Range("B76").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C76").Value = "Arica y Parinacota"
$ws.Range("D76").Value = 45236
$ws.Range("E76").Value = 15
$ws.Range("F76").Value = 100112038
$ws.Range("G76").Value = "Cebollín baby"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 300
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = 3250
$ws.Range("N76").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O76").Value = "Región de Arica y Parinacota"
$ws.Range("P76").Value = 1625
$ws.Range("Q76").Value = 2
$ws.Range("R76").Value = "Hortaliza"
